$d = $word.ActiveDocument

# The document embeds the Pearson logo twice (in the two footers) and the
# BTec logo once (in the header). Each inline picture carries its
# "insertion name" in two places in the drawing markup: <wp:docPr name="..."/>
# and <pic:cNvPr name="..."/>. Word renumbered these names:
#   Pearson logo: image1.png -> image2.png  (both footers)
#   BTec logo:    image2.jpg -> image1.jpg  (header)
# The alt text / description (descr="...") and the actual embedded media
# parts are untouched - only these cosmetic name= attributes change.

$xml = $d.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.WordOpenXML = $xml
